$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# --- Crime data table updates (rows 14-27) ---
$ws.Range("F30").Copy($ws.Range("C14"))
$ws.Range("C14").Value = 1
$ws.Range("F30").Copy($ws.Range("F14"))
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("F30").Copy($ws.Range("I14"))
$ws.Range("I14").Value = 1
$ws.Range("K14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("C23").Copy($ws.Range("C15"))
$ws.Range("F30").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 2
$ws.Range("K37").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("F30").Copy($ws.Range("G15"))
$ws.Range("G15").Value = 2
$ws.Range("K37").Copy($ws.Range("H15"))
$ws.Range("H15").Value = -50
$ws.Range("F30").Copy($ws.Range("J15"))
$ws.Range("J15").Value = 2
$ws.Range("K37").Copy($ws.Range("K15"))
$ws.Range("K15").Value = -50
$ws.Range("L15").Value = -66.666666666666
$ws.Range("F30").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 1
$ws.Range("K37").Copy($ws.Range("E16"))
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 2
$ws.Range("F30").Copy($ws.Range("G16"))
$ws.Range("G16").Value = 1
$ws.Range("K37").Copy($ws.Range("H16"))
$ws.Range("H16").Value = 100
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 266.666666666667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -26.666666666666
$ws.Range("N16").Value = -86.746987951807
$ws.Range("C17").Value = 8
$ws.Range("F30").Copy($ws.Range("D17"))
$ws.Range("D17").Value = 1
$ws.Range("K37").Copy($ws.Range("E17"))
$ws.Range("E17").Value = 700
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = 300
$ws.Range("I17").Value = 25
$ws.Range("J17").Value = 8
$ws.Range("K17").Value = 212.5
$ws.Range("L17").Value = 127.272727272727
$ws.Range("M17").Value = 108.333333333333
$ws.Range("N17").Value = -35.897435897435
$ws.Range("F30").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 77.777777777777
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 13
$ws.Range("K18").Value = 100
$ws.Range("L18").Value = 116.666666666667
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = -84.33734939759
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 5.882352941176
$ws.Range("I19").Value = 53
$ws.Range("J19").Value = 75
$ws.Range("K19").Value = -29.333333333333
$ws.Range("L19").Value = 47.222222222222
$ws.Range("M19").Value = 17.777777777777
$ws.Range("N19").Value = -31.168831168831
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 12
$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 16
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = 23.076923076923
$ws.Range("L20").Value = 166.666666666667
$ws.Range("M20").Value = -36
$ws.Range("N20").Value = -94.444444444444
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 53.333333333333
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = 52.727272727272
$ws.Range("I21").Value = 133
$ws.Range("J21").Value = 115
$ws.Range("K21").Value = 15.652173913043
$ws.Range("L21").Value = 68.354430379746
$ws.Range("M21").Value = 13.675213675213
$ws.Range("N21").Value = -79.817905918057
$ws.Range("F30").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K37").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("F30").Copy($ws.Range("G22"))
$ws.Range("G22").Value = 1
$ws.Range("K37").Copy($ws.Range("H22"))
$ws.Range("H22").Value = -100
$ws.Range("F30").Copy($ws.Range("J22"))
$ws.Range("J22").Value = 1
$ws.Range("K37").Copy($ws.Range("K22"))
$ws.Range("K22").Value = -100
$ws.Range("K37").Copy($ws.Range("L22"))
$ws.Range("L22").Value = -100
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -13.513513513513
$ws.Range("F24").Value = 119
$ws.Range("G24").Value = 156
$ws.Range("H24").Value = -23.717948717948
$ws.Range("I24").Value = 239
$ws.Range("J24").Value = 264
$ws.Range("K24").Value = -9.469696969696
$ws.Range("L24").Value = 97.520661157024
$ws.Range("M24").Value = 36.571428571428
$ws.Range("C25").Value = 6
$ws.Range("E25").Value = -14.285714285714
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = 8
$ws.Range("I25").Value = 59
$ws.Range("J25").Value = 42
$ws.Range("K25").Value = 40.47619047619
$ws.Range("L25").Value = 78.787878787878
$ws.Range("M25").Value = 13.461538461538
$ws.Range("C23").Copy($ws.Range("C26"))
$ws.Range("F30").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 2
$ws.Range("K37").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
$ws.Range("F30").Copy($ws.Range("G26"))
$ws.Range("G26").Value = 2
$ws.Range("K37").Copy($ws.Range("H26"))
$ws.Range("H26").Value = -50
$ws.Range("F30").Copy($ws.Range("J26"))
$ws.Range("J26").Value = 2
$ws.Range("K37").Copy($ws.Range("K26"))
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = -33.333333333333
$ws.Range("C23").Copy($ws.Range("C27"))
$ws.Range("E27").Value = -100
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = -80
$ws.Range("L27").Value = -66.666666666666

Write-Host "Edit complete"